{"js": "// Popular-science abstract (Danish) proofreading pass:\n//   - \"patienterne\" -> \"patienter\"\n//   - \"sikre a blot\" -> \"sikre at blot\" (missing \"t\" typo fix)\n//   - \"modeller til at forudsige\" -> \"modeller for at forudsige\"\n//   - \"af nye molekyler.\" -> \"af nye potentielle absorptionsfremmere.\"\n//   - wrap the word \"absorptionsfremmere\" (in \"som absorptionsfremmere\") with a bookmark\n//   - \"forskere p\u00e5 tv\u00e6rs\" -> \"forskning p\u00e5 tv\u00e6rs\"\n\nconst body = context.document.body;\n\n// 1) \"...for patienterne.\" -> \"...for patienter.\"\nlet found = body.search(\"patienterne\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"patienter\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"...sikre a blot...\" -> \"...sikre at blot...\"\nfound = body.search(\"sikre a blot\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"sikre at blot\", \"Replace\");\n  await context.sync();\n}\n\n// 3) \"modeller til at forudsige\" -> \"modeller for at forudsige\"\nfound = body.search(\"modeller til at forudsige\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"modeller for at forudsige\", \"Replace\");\n  await context.sync();\n}\n\n// 4) \"af nye molekyler. Modeller af random forest\" -> \"af nye potentielle absorptionsfremmere. Modeller af random forest\"\nfound = body.search(\"af nye molekyler. Modeller af random forest\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\n    \"af nye potentielle absorptionsfremmere. Modeller af random forest\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 5) Wrap \"absorptionsfremmere\" in \"som absorptionsfremmere\" with the __DdeLink__ bookmark\nfound = body.search(\"som absorptionsfremmere\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (found.items.length > 0) {\n  const word = found.items[0].search(\"absorptionsfremmere\", { matchCase: true, matchWholeWord: false });\n  await context.sync();\n  if (word.items.length > 0) {\n    word.items[0].insertBookmark(\"__DdeLink__4265_2146180382\");\n    await context.sync();\n  }\n}\n\n// 6) \"forskere p\u00e5 tv\u00e6rs\" -> \"forskning p\u00e5 tv\u00e6rs\"\nfound = body.search(\"forskere p\u00e5 tv\u00e6rs\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"forskning p\u00e5 tv\u00e6rs\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Popular-science abstract (Danish) proofreading pass:\n#   - \"patienterne\" -> \"patienter\"\n#   - \"sikre a blot\" -> \"sikre at blot\" (missing \"t\" typo fix)\n#   - \"modeller til at forudsige\" -> \"modeller for at forudsige\"\n#   - \"af nye molekyler.\" -> \"af nye potentielle absorptionsfremmere.\"\n#   - wrap the word \"absorptionsfremmere\" (in \"som absorptionsfremmere\") with a bookmark\n#   - \"forskere p\u00e5 tv\u00e6rs\" -> \"forskning p\u00e5 tv\u00e6rs\"\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2, wdFindContinue = 1 (literal values; Word enums aren't predefined here)\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceAll)\n}\n\n# 1) \"...for patienterne.\" -> \"...for patienter.\"\nReplace-Text \"patienterne\" \"patienter\"\n\n# 2) \"...sikre a blot...\" -> \"...sikre at blot...\"\nReplace-Text \"sikre a blot\" \"sikre at blot\"\n\n# 3) \"modeller til at forudsige\" -> \"modeller for at forudsige\"\nReplace-Text \"modeller til at forudsige\" \"modeller for at forudsige\"\n\n# 4) \"af nye molekyler. Modeller af random forest\" -> \"af nye potentielle absorptionsfremmere. Modeller af random forest\"\nReplace-Text \"af nye molekyler. Modeller af random forest\" \"af nye potentielle absorptionsfremmere. Modeller af random forest\"\n\n# 5) Wrap \"absorptionsfremmere\" in \"som absorptionsfremmere\" with the __DdeLink__ bookmark\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"som absorptionsfremmere\"\n$find.Execute()\n$wordRange = $find.Parent\n$wordRange.Start = $wordRange.Start + 4\n$d.Bookmarks.Add(\"__DdeLink__4265_2146180382\", $wordRange)\n\n# 6) \"forskere p\u00e5 tv\u00e6rs\" -> \"forskning p\u00e5 tv\u00e6rs\"\nReplace-Text \"forskere p\u00e5 tv\u00e6rs\" \"forskning p\u00e5 tv\u00e6rs\"\n"}
